# Updating Trans specs and Frag 1 a

$wb = $excel.ActiveWorkbook
$charsWs = $wb.Worksheets.Item("CHARs")

# --- Update data validation lists on CHARs (rows 1-7) ---
# I: add "erased" option
# K, L, N: add "null" option
# M: add "null" option
# O: add "null" option
$iFormula = '"transformed,reinked,retraced,reinked?,retraced?,intralinear,creased,erased"'
$boolFormula = '"null,True,False"'
$mFormula = '"null,False,True,relevant_w,relevant_h"'
$oFormula = '"null,certain,probable_letter,possible_letter"'

for ($row = 1; $row -le 7; $row++) {
    $charsWs.Range("I$row").Validation.Formula1 = $iFormula
    $charsWs.Range("K$row").Validation.Formula1 = $boolFormula
    $charsWs.Range("L$row").Validation.Formula1 = $boolFormula
    $charsWs.Range("M$row").Validation.Formula1 = $mFormula
    $charsWs.Range("N$row").Validation.Formula1 = $boolFormula
    $charsWs.Range("O$row").Validation.Formula1 = $oFormula
}

# --- Add new "Sub_Frags" worksheet after the existing sheets ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$subFrags = $wb.Worksheets.Add($null, $lastSheet)
$subFrags.Name = "Sub_Frags"

$headers = @("frag_id", "iaa_img_id", "Label", "Area", "Mean", "Min", "Max", "BX", "BY", "Width", "Height", "Major", "Minor", "Circ.", "AR", "Round", "Solidity")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $subFrags.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Keep the original active sheet selected (adding a sheet makes it active by default)
$charsWs.Activate()
